$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns keep their text representation (as in the source
# workbook, where these values are stored as inline strings, not numbers).
$ws.Range("D2:E47").NumberFormat = "@"

$ws.Range("D2").Value = "286.32"
$ws.Range("E2").Value = "3.85%"
$ws.Range("D3").Value = "28.45"
$ws.Range("E3").Value = "4.69%"
$ws.Range("E4").Value = "1.49%"
$ws.Range("D5").Value = "0.06537"
$ws.Range("E5").Value = "2.24%"
$ws.Range("D6").Value = "7.249"
$ws.Range("E6").Value = "4.30%"
$ws.Range("D7").Value = "1.374"
$ws.Range("E7").Value = "16.09%"
$ws.Range("D8").Value = "0.9104"
$ws.Range("E8").Value = "3.99%"
$ws.Range("D9").Value = "0.1563"
$ws.Range("E9").Value = "3.55%"
$ws.Range("D10").Value = "0.06758"
$ws.Range("E10").Value = "33.16%"
$ws.Range("D11").Value = "0.07723"
$ws.Range("E11").Value = "1.91%"
$ws.Range("D12").Value = "0.02975"
$ws.Range("E12").Value = "0.64%"
$ws.Range("D13").Value = "0.08967"
$ws.Range("E13").Value = "-0.17%"
$ws.Range("D14").Value = "0.001604"
$ws.Range("E14").Value = "2.36%"
$ws.Range("D15").Value = "0.0006540"
$ws.Range("E15").Value = "2.39%"
$ws.Range("D16").Value = "0.006012"
$ws.Range("E16").Value = "-2.74%"
$ws.Range("D17").Value = "3.468"
$ws.Range("E17").Value = "-0.20%"
$ws.Range("D18").Value = "3.392"
$ws.Range("E18").Value = "2.58%"
$ws.Range("D19").Value = "2.236"
$ws.Range("E19").Value = "-0.73%"
$ws.Range("D20").Value = "0.3150"
$ws.Range("E20").Value = "0.46%"
$ws.Range("E21").Value = "0.10%"
$ws.Range("D22").Value = "4.022"
$ws.Range("E22").Value = "2.87%"
$ws.Range("D23").Value = "0.1554"
$ws.Range("E23").Value = "12.63%"
$ws.Range("D24").Value = "0.04460"
$ws.Range("E24").Value = "1.10%"
$ws.Range("D25").Value = "0.001191"
$ws.Range("E25").Value = "0.96%"
$ws.Range("D26").Value = "0.004332"
$ws.Range("E26").Value = "12.27%"
$ws.Range("E28").Value = "-1.70%"
$ws.Range("D29").Value = "0.0001636"
$ws.Range("E29").Value = "-15.71%"
$ws.Range("D40").Value = "0.04165"
$ws.Range("E40").Value = "0.29%"
$ws.Range("D41").Value = "0.1416"
$ws.Range("E41").Value = "20.69%"
$ws.Range("D42").Value = "0.006637"
$ws.Range("E42").Value = "-3.39%"
$ws.Range("D43").Value = "0.002160"
$ws.Range("E43").Value = "-1.40%"
$ws.Range("D44").Value = "0.01241"
$ws.Range("E44").Value = "5.01%"
$ws.Range("D45").Value = "0.00005571"
$ws.Range("E45").Value = "7.08%"
$ws.Range("D46").Value = "1.562"
$ws.Range("E46").Value = "-6.93%"
$ws.Range("D47").Value = "0.01850"
